# Weekly "cryptos" snapshot refresh (GitHub Actions scraper run).
# Updates prices/volumes, rotates several coin rows that changed rank
# order, and bumps the "Hora" (hour) column from 2 to 3 for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, applied in one pass below.
$updates = [ordered]@{
    "D2" = "240.50"
    "G2" = "3"
    "G3" = "3"
    "D4" = "5.535"
    "G4" = "3"
    "D5" = "0.05587"
    "G5" = "3"
    "D6" = "3.380"
    "G6" = "3"
    "D7" = "6.471"
    "G7" = "3"
    "D8" = "1.088"
    "G8" = "3"
    "D9" = "0.8012"
    "G9" = "3"
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D10" = "0.1421"
    "E10" = "9WazirXWRX"
    "G10" = "3"
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D11" = "0.07401"
    "E11" = "10MandalaExchangeTokenMDX"
    "G11" = "3"
    "B12" = "LiechtensteinCryptoassetsExchange"
    "C12" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D12" = "0.03243"
    "E12" = "11LiechtensteinCryptoassetsExchangeLCX"
    "G12" = "3"
    "B13" = "BitrueCoin"
    "C13" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D13" = "0.02992"
    "E13" = "12BitrueCoinBTR"
    "G13" = "3"
    "B14" = "BitMartToken"
    "C14" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D14" = "0.09252"
    "E14" = "13BitMartTokenBMX"
    "G14" = "3"
    "B15" = "BitForexToken"
    "C15" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D15" = "0.001661"
    "E15" = "14BitForexTokenBF"
    "G15" = "3"
    "B16" = "MCDex"
    "C16" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "D16" = "3.266"
    "E16" = "15MCDexMCB"
    "G16" = "3"
    "B17" = "CoinExToken"
    "C17" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "D17" = "0.04707"
    "E17" = "16CoinExTokenCET"
    "G17" = "3"
    "B18" = "One"
    "C18" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "D18" = "0.0005739"
    "E18" = "17OneONE"
    "G18" = "3"
    "D19" = "0.006260"
    "G19" = "3"
    "G20" = "3"
    "G21" = "3"
    "G22" = "3"
    "G23" = "3"
    "D24" = "3.978"
    "G24" = "3"
    "D25" = "2.129"
    "G25" = "3"
    "G26" = "3"
    "G27" = "3"
    "G28" = "3"
    "G29" = "3"
    "G30" = "3"
    "G31" = "3"
    "G32" = "3"
    "G33" = "3"
    "G34" = "3"
    "G35" = "3"
    "G36" = "3"
    "G37" = "3"
    "G38" = "3"
    "G39" = "3"
    "D40" = "0.04179"
    "E40" = "39IDEXIDEXBestin24h"
    "G40" = "3"
    "B41" = "BKEXToken"
    "C41" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "D41" = "0.1045"
    "E41" = "40BKEXTokenBKK"
    "G41" = "3"
    "D42" = "0.002970"
    "G42" = "3"
    "B43" = "KickToken"
    "C43" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "D43" = "0.003254"
    "E43" = "42KickTokenKICK"
    "G43" = "3"
    "D44" = "0.009158"
    "G44" = "3"
    "D45" = "0.00005504"
    "G45" = "3"
    "G46" = "3"
    "G47" = "3"
    "D48" = "0.03059"
    "G48" = "3"
    "G49" = "3"
    "G50" = "3"
    "G51" = "3"
}

foreach ($addr in $updates.Keys) {
    # Keep these cells stored as text (not auto-converted to numbers),
    # matching the source data which preserves exact formatting such as
    # trailing/leading zeros (e.g. "240.50", "0.00005504").
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
